$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2823
$ws.Range("E2").Value = 103
$ws.Range("F2").Value = 103
$ws.Range("G2").Value = 192
$ws.Range("H2").Value = 154
$ws.Range("I2").Value = 154
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2440
$ws.Range("L2").Value = 715
$ws.Range("M2").Value = 1725
$ws.Range("N2").Value = 1725
$ws.Range("O2").Value = 122
$ws.Range("P2").Value = 67
$ws.Range("Q2").Value = 86
$ws.Range("R2").Value = -137
$ws.Range("S2").Value = 110
$ws.Range("T2").Value = 5
$ws.Range("U2").Value = 80
$ws.Range("V2").Value = 104
$ws.Range("W2").Value = 3.66
$ws.Range("X2").Value = 5.45
$ws.Range("Y2").Value = 9.640000000000001
$ws.Range("Z2").Value = 6.65
$ws.Range("AA2").Value = 41.45
$ws.Range("AB2").Value = 2466.97
$ws.Range("AC2").Value = 571
$ws.Range("AD2").Value = 8.550000000000001
$ws.Range("AE2").Value = 6406
$ws.Range("AF2").Value = 0.76
$ws.Range("AG2").Value = 375
$ws.Range("AH2").Value = 7.68
$ws.Range("AI2").Value = 65.65000000000001
$ws.Range("AJ2").Value = 26958990

# Row 3
$ws.Range("D3").Value = 2960
$ws.Range("E3").Value = 90
$ws.Range("F3").Value = 90
$ws.Range("G3").Value = 25
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 29
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2444
$ws.Range("L3").Value = 750
$ws.Range("M3").Value = 1694
$ws.Range("N3").Value = 1693
$ws.Range("O3").Value = 181
$ws.Range("P3").Value = 135
$ws.Range("Q3").Value = 378
$ws.Range("R3").Value = -117
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 375
$ws.Range("V3").Value = 165
$ws.Range("W3").Value = 3.04
$ws.Range("X3").Value = 0.97
$ws.Range("Y3").Value = 1.72
$ws.Range("Z3").Value = 1.18
$ws.Range("AA3").Value = 44.29
$ws.Range("AB3").Value = 1132.2
$ws.Range("AC3").Value = 109
$ws.Range("AD3").Value = 34.87
$ws.Range("AE3").Value = 6375
$ws.Range("AF3").Value = 0.6
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 5.27
$ws.Range("AI3").Value = 180.92
$ws.Range("AJ3").Value = 26999207

# Row 4
$ws.Range("D4").Value = 2636
$ws.Range("E4").Value = 78
$ws.Range("F4").Value = 78
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = -1
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = -15
$ws.Range("K4").Value = 2844
$ws.Range("L4").Value = 1018
$ws.Range("M4").Value = 1826
$ws.Range("N4").Value = 1787
$ws.Range("O4").Value = 230
$ws.Range("P4").Value = 135
$ws.Range("Q4").Value = -66
$ws.Range("R4").Value = -153
$ws.Range("S4").Value = -11
$ws.Range("T4").Value = 2
$ws.Range("U4").Value = -68
$ws.Range("V4").Value = 243
$ws.Range("W4").Value = 2.97
$ws.Range("X4").Value = -0.02
$ws.Range("Y4").Value = 0.02
$ws.Range("Z4").Value = -0.02
$ws.Range("AA4").Value = 55.76
$ws.Range("AB4").Value = 1098.85
$ws.Range("AC4").Value = 1
$ws.Range("AD4").Value = 7056.13
$ws.Range("AE4").Value = 6718
$ws.Range("AF4").Value = 1.09
$ws.Range("AG4").Value = 80
$ws.Range("AH4").Value = 1.09
$ws.Range("AI4").Value = 7585.26
$ws.Range("AJ4").Value = 27039689

# Row 5
$ws.Range("D5").Value = 2461
$ws.Range("E5").Value = 88
$ws.Range("F5").Value = 88
$ws.Range("G5").Value = 133
$ws.Range("H5").Value = 88
$ws.Range("I5").Value = 98
$ws.Range("J5").Value = 15
$ws.Range("K5").Value = 2669
$ws.Range("L5").Value = 863
$ws.Range("M5").Value = 1807
$ws.Range("N5").Value = 1777
$ws.Range("O5").Value = 251
$ws.Range("P5").Value = 135
$ws.Range("Q5").Value = 21
$ws.Range("R5").Value = -13
$ws.Range("S5").Value = -58
$ws.Range("T5").Value = 13
$ws.Range("U5").Value = 8
$ws.Range("V5").Value = 236
$ws.Range("W5").Value = 3.57
$ws.Range("X5").Value = 3.56
$ws.Range("Y5").Value = 5.48
$ws.Range("Z5").Value = 3.18
$ws.Range("AA5").Value = 47.74
$ws.Range("AB5").Value = 1154.4
$ws.Range("AC5").Value = 361
$ws.Range("AD5").Value = 26.75
$ws.Range("AE5").Value = 6677
$ws.Range("AF5").Value = 1.45
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.04
$ws.Range("AI5").Value = 27.25
$ws.Range("AJ5").Value = 27049809

# Row 6
$ws.Range("D6").Value = 2703
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 166
$ws.Range("H6").Value = 106
$ws.Range("I6").Value = 120
$ws.Range("K6").Value = 3241
$ws.Range("L6").Value = 1325
$ws.Range("M6").Value = 1916
$ws.Range("N6").Value = 1881
$ws.Range("P6").Value = 136
$ws.Range("Q6").Value = 90
$ws.Range("R6").Value = 52
$ws.Range("S6").Value = -170
$ws.Range("T6").Value = 11
$ws.Range("U6").Value = 79
$ws.Range("V6").Value = 163
$ws.Range("W6").Value = 0.76
$ws.Range("X6").Value = 3.94
$ws.Range("Y6").Value = 6.55
$ws.Range("Z6").Value = 3.6
$ws.Range("AA6").Value = 69.16
$ws.Range("AB6").Value = 1307.26
$ws.Range("AC6").Value = 442
$ws.Range("AD6").Value = 14.45
$ws.Range("AE6").Value = 7029
$ws.Range("AF6").Value = 0.91
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 1.57
$ws.Range("AI6").Value = 22.33
$ws.Range("AJ6").Value = 27190652

# Row 7: clear then set only surviving cells
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D7").Value = 3211
$ws.Range("E7").Value = 88
$ws.Range("G7").Value = 283
$ws.Range("H7").Value = 250
$ws.Range("W7").Value = 2.74
$ws.Range("X7").Value = 7.79

# Rows 8 and 9: clear all data cells, keep only A,B,C
$ws.Range("D8:AJ9").ClearContents()
